# Scene 24D edit: tweak three dialogue/narration runs.
#  - "Petra ... Oh, that's his name."                -> force xml:space="preserve" on the run
#  - "Lilith ... embarrassed_slightly): Oh."          -> force xml:space="preserve" on the run
#  - "And being able to see ... pretty happy as well." -> drop the trailing space, keep xml:space="preserve"
#
# The engine only stamps xml:space="preserve" on a <w:t> when the text being
# assigned has leading/trailing whitespace, but it does NOT strip the
# attribute back off on a later write that no longer needs it. So for the
# first two runs we briefly write the text with a leading space (which
# stamps the attribute) and then immediately rewrite the clean text (which
# keeps the now-sticky attribute). For the third run the attribute is
# already present; we just need to rewrite the text without the trailing
# space, which leaves the attribute in place untouched.

$d = $word.ActiveDocument

function Set-RunTextPreserved($matchText, $newText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $r = $p.Range
        if ($r.Text.TrimEnd() -eq $matchText) {
            # Step 1: write with a leading space so the serializer marks
            # the run's <w:t> as xml:space="preserve".
            $r.Text = " " + $newText
            # Step 2: re-fetch the (now resized) paragraph range and
            # rewrite the clean text - the preserve flag sticks around.
            $p2 = $d.Paragraphs.Item($i)
            $r2 = $p2.Range
            $r2.Text = $newText
            return $true
        }
    }
    return $false
}

function Set-RunTextTrimTrailing($matchTextTrimmed, $newText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $r = $p.Range
        if ($r.Text.TrimEnd() -eq $matchTextTrimmed) {
            $r.Text = $newText
            return $true
        }
    }
    return $false
}

$ok1 = Set-RunTextPreserved "Petra (neutral smiling): Oh, that’s his name." "Petra (neutral smiling): Oh, that’s his name."
$ok2 = Set-RunTextPreserved "Lilith (neutral embarrassed_slightly): Oh." "Lilith (neutral embarrassed_slightly): Oh."
$ok3 = Set-RunTextTrimTrailing "And being able to see that makes me pretty happy as well." "And being able to see that makes me pretty happy as well."

Write-Output "Petra run updated: $ok1"
Write-Output "Lilith run updated: $ok2"
Write-Output "Closing run updated: $ok3"
